$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The trailing "separator" rows (16, 31, 46, 53) were blank spacer rows
# (only J/K carried the numeric-cell style). Remove them so the data
# blocks below each one shift up and become contiguous - this matches
# the final filtered/cleaned layout (dimension shrinks from K55 to K51).
# Deleting from the bottom up keeps the remaining row numbers stable.
$ws.Rows.Item(53).Delete() | Out-Null
$ws.Rows.Item(46).Delete() | Out-Null
$ws.Rows.Item(31).Delete() | Out-Null
$ws.Rows.Item(16).Delete() | Out-Null

# Restore the cursor/selection to where the author left it after editing.
$ws.Range("D34").Select() | Out-Null
